# Fixed all Dialog GUIs except Stage Related. Added Aux Classes
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Stages sheet: wire H5:H20 up to formulas pulling from "Participants @Events"
#    and give the two date columns (B, J) the existing short-date formatting.
# ---------------------------------------------------------------------------
$stages = $wb.Worksheets.Item("Stages")

for ($r = 5; $r -le 20; $r++) {
    $stages.Range("H$r").Formula = "='Participants @Events'!E$r"
}

# Reuse the workbook's existing date number format (already used on sheet
# "Events, Teams, Equipment" D4) instead of typing a format string, so Excel
# reuses the existing style record rather than fabricating a new one.
$dateSample = $wb.Worksheets.Item("Events, Teams, Equipment").Range("D4")

$dateSample.Copy()
$stages.Range("B5:B20").PasteSpecial(-4122)   # xlPasteFormats

$dateSample.Copy()
$stages.Range("J5:J20").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

$stages.Activate() | Out-Null
$stages.Range("H5:M20").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2) Add a new "Sheet1" with small auxiliary lookup tables (DOCUMENT,
#    OTHERFILE, TEXTFILE) used by the "Aux Classes".
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$aux = $wb.Worksheets.Add($null, $lastSheet)
$aux.Name = "Sheet1"

# DOCUMENT table (id, eventNumber)
$aux.Range("B3").Value = "DOCUMENT"
$aux.Range("B4").Value = "id"
$aux.Range("C4").Value = "eventNumber"
$aux.Range("B5").Value = 1
$aux.Range("C5").Value = 1
$aux.Range("B6").Value = 2
$aux.Range("C6").Value = 1
$aux.Range("B7").Value = 3
$aux.Range("C7").Value = 2
$aux.Range("B8").Value = 4
$aux.Range("C8").Value = 1

# OTHERFILE table (id, path)
$aux.Range("F3").Value = "OTHERFILE"
$aux.Range("F4").Value = "id"
$aux.Range("G4").Value = "path"
$aux.Range("F5").Value = 1
$aux.Range("G5").Value = "MyPath1"
$aux.Range("F6").Value = 2
$aux.Range("G6").Value = "MyPath2"

# TEXTFILE table (id, content)
$aux.Range("I3").Value = "TEXTFILE"
$aux.Range("I4").Value = "id"
$aux.Range("J4").Value = "content"
$aux.Range("I5").Value = 3
$aux.Range("J5").Value = "Test"
$aux.Range("I6").Value = 4
$aux.Range("J6").Value = "Test2"

$aux.Activate() | Out-Null
$aux.Range("I4:J6").Select() | Out-Null
